# Update Name of Algo
# Applies updated imputed values to the RandomForest result data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value  = 5.608099999999991
$ws.Range("E6").Value  = 12.4394
$ws.Range("E7").Value  = 11.9996
$ws.Range("B8").Value  = 5.326899999999991
$ws.Range("E8").Value  = 13.9124
$ws.Range("A12").Value = -22.77060000000001
$ws.Range("B12").Value = 5.953100000000001
$ws.Range("B14").Value = 8.823000000000006
$ws.Range("E19").Value = 13.0703
$ws.Range("E21").Value = 12.66899999999999
$ws.Range("B22").Value = 5.008400000000004
$ws.Range("E24").Value = 12.75189999999999

$wb.Save()
